$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Re-home the existing header text one column to the right (B->C->D->E)
#    so the brand new "Computer" column can be inserted at B. Order among
#    these three doesn't matter for correctness (final state only), but we
#    do it before introducing any new strings so nothing is lost.
# ---------------------------------------------------------------------------
$ws.Range("C1").Value = "Function"
$ws.Range("D1").Value = "n tips"
$ws.Range("E1").Value = "n timesteps"

# ---------------------------------------------------------------------------
# 2. Re-home the existing "Sim" / "Results" function values from column B
#    into column C for every data row that needs them.
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "Sim"
$ws.Range("C3").Value = "Results"
$ws.Range("C4").Value = "Sim"
$ws.Range("C5").Value = "Sim"
$ws.Range("C6").Value = "Sim"

# ---------------------------------------------------------------------------
# 3. Introduce the brand-new strings. The order below fixes the order they
#    land in the shared-string table, so keep it exactly as-is.
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "Walltime per core (hours)"
$ws.Range("F1").Value = "calculation steps (tips * time)"
$ws.Range("H1").Value = "Projected time to run (hours)"
$ws.Range("B1").Value = "Computer"
$ws.Range("B2").Value = "8-core Macpro"
$ws.Range("B3").Value = "6-core Macpro"
$ws.Range("B4").Value = "6-core Macpro"
$ws.Range("G4").Value = "4.8 (ran as if 8 cores, so this number is misleading)"

# ---------------------------------------------------------------------------
# 4. Fill in the rest of the numeric grid (existing rows 2-3 keep their n
#    tips / n timesteps / walltime numbers; new rows 4-6 are today's runs).
# ---------------------------------------------------------------------------
$ws.Range("D2").Value = 300
$ws.Range("E2").Value = 300
$ws.Range("G2").Value = 3.57

$ws.Range("D3").Value = 300
$ws.Range("E3").Value = 300
$ws.Range("G3").Value = 8

$ws.Range("A4").Value = $ws.Range("A2").Value2
$ws.Range("D4").Value = 300
$ws.Range("E4").Value = 300

$ws.Range("A5").Value = 42572
$ws.Range("D5").Value = 1200
$ws.Range("E5").Value = 300

$ws.Range("A6").Value = 42572
$ws.Range("D6").Value = 1200
$ws.Range("E6").Value = 300

# ---------------------------------------------------------------------------
# 5. Date formatting for column A (already in place for rows 2-3) and for
#    the new column B "Computer" cells, which inherited the date format too.
# ---------------------------------------------------------------------------
$ws.Range("A4:A6").NumberFormat = "d-mmm-yy"
$ws.Range("B2:B6").NumberFormat = "d-mmm-yy"

# ---------------------------------------------------------------------------
# 6. Formulas: F2 is a standalone formula, F3:F24 share one formula.
# ---------------------------------------------------------------------------
$ws.Range("F2").Formula = "=D2*E2"
$ws.Range("F3:F24").Formula = "=D3*E3"

$ws.Range("H5").Formula = "=(F5/F2)*G`$2"
$ws.Range("H6").Formula = "=(F6/F3)*G2"

# ---------------------------------------------------------------------------
# 7. Right-align the "Walltime per core" / "Projected time" columns.
# ---------------------------------------------------------------------------
$ws.Range("G1:G4").HorizontalAlignment = -4152
$ws.Range("H1").HorizontalAlignment = -4152
$ws.Range("H5:H6").HorizontalAlignment = -4152

# ---------------------------------------------------------------------------
# 8. Column widths to match the new layout.
# ---------------------------------------------------------------------------
$ws.Columns("F").ColumnWidth = 24.666666666666668
$ws.Columns("G").ColumnWidth = 42.498697916666664
$ws.Columns("H").ColumnWidth = 24.166666666666668

# ---------------------------------------------------------------------------
# 9. Selection cursor ends on F9, matching the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("F9").Select()
